$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 13.36538210713172
    "C2" = 6.695632360840241
    "D2" = 7.963456263971554
    "E2" = 12.84960670177668
    "F2" = 37.25589866466917
    "I2" = 28.00070937863699
    "J2" = 10.14632020270577
    "K2" = 10.74762511573998
    "L2" = 11.19551048030369
    "M2" = 15.38524925947976
    "N2" = 21.41455861298659
    "O2" = 28.77495283335093
    "B3" = 13.19134087224952
    "C3" = 6.634402995201029
    "D3" = 7.954597282926569
    "E3" = 12.87095082194856
    "F3" = 37.34020316456343
    "I3" = 28.08630121466519
    "J3" = 10.16130629157166
    "K3" = 10.62279215314741
    "L3" = 11.20415672931712
    "M3" = 15.36437822379487
    "N3" = 21.47173500941641
    "O3" = 28.85430788210193
    "B4" = 13.08553309970777
    "C4" = 6.595948487652365
    "D4" = 7.950142678483177
    "E4" = 12.88524882360572
    "F4" = 37.39844483412981
    "I4" = 28.14282860660882
    "J4" = 10.17100951227416
    "K4" = 10.54703833805142
    "L4" = 11.21062491242877
    "M4" = 15.35341062579724
    "N4" = 21.50851003005478
    "O4" = 28.90756208757643
    "B5" = 13.04272982242598
    "C5" = 6.580068127104761
    "D5" = 7.948576769787878
    "E5" = 12.89137582664983
    "F5" = 37.42380704983395
    "I5" = 28.16686369523543
    "J5" = 10.17509016190037
    "K5" = 10.51642546351323
    "L5" = 11.21355297384081
    "M5" = 15.34940944668438
    "N5" = 21.5239168964715
    "O5" = 28.93040249528845
    "B6" = 13.035642779749
    "C6" = 6.577418716600708
    "D6" = 7.948331862024081
    "E6" = 12.89241137357041
    "F6" = 37.4281167434979
    "I6" = 28.17091509896867
    "J6" = 10.17577540180916
    "K6" = 10.51135871092006
    "L6" = 11.21405684551831
    "M6" = 15.34877343558607
    "N6" = 21.52650064211404
    "O6" = 28.93426390737401
    "B7" = 13.08495450269525
    "C7" = 6.595735160830898
    "D7" = 7.95012054815741
    "E7" = 12.88533023725595
    "F7" = 37.39878028621671
    "I7" = 28.14314870372121
    "J7" = 10.17106403264551
    "K7" = 10.54662439708299
    "L7" = 11.21066321718544
    "M7" = 15.35335476404914
    "N7" = 21.50871610727061
    "O7" = 28.90786551027027
    "B8" = 13.30518345375276
    "C8" = 6.674702583014745
    "D8" = 7.960198353692187
    "E8" = 12.85671896027803
    "F8" = 37.28362185379812
    "I8" = 28.02939701832913
    "J8" = 10.15138348991171
    "K8" = 10.70441678833007
    "L8" = 11.19825151000534
    "M8" = 15.37767179600848
    "N8" = 21.43392746270413
    "O8" = 28.80137417392836
    "B9" = 13.74324987089821
    "C9" = 6.822469423968106
    "D9" = 7.987695149408634
    "E9" = 12.81005065813148
    "F9" = 37.10922935887071
    "I9" = 27.83783657662017
    "J9" = 10.11675472726563
    "K9" = 11.01949735549784
    "L9" = 11.18307993518452
    "M9" = 15.43984487921818
    "N9" = 21.30045042985808
    "O9" = 28.62849532349913
    "B10" = 14.06598710138397
    "C10" = 6.926362933717598
    "D10" = 8.012498875185969
    "E10" = 12.78148389353953
    "F10" = 37.01249160027564
    "I10" = 27.71627246617028
    "J10" = 10.09370748013383
    "K10" = 11.25249355141858
    "L10" = 11.17747944793907
    "M10" = 15.49411566023979
    "N10" = 21.21034289319119
    "O10" = 28.52340685459263
    "B11" = 14.21240084361255
    "C11" = 6.972540807652806
    "D11" = 8.024755523230342
    "E11" = 12.76972308288734
    "F11" = 36.97530221019591
    "I11" = 27.66512701488303
    "J11" = 10.08373785829158
    "K11" = 11.35840922530282
    "L11" = 11.17612608431001
    "M11" = 15.52061483275132
    "N11" = 21.17106183898338
    "O11" = 28.48036166203184
    "B12" = 14.26773739768693
    "C12" = 6.989865389499955
    "D12" = 8.029534310908137
    "E12" = 12.76544647895781
    "F12" = 36.96219975117225
    "I12" = 27.64635657077945
    "J12" = 10.08003627420844
    "K12" = 11.39847242479214
    "L12" = 11.17578443414196
    "M12" = 15.53090469009548
    "N12" = 21.15643168496026
    "O12" = 28.46474614007378
    "B13" = 14.2558252192155
    "C13" = 6.98614153066887
    "D13" = 8.028499040854225
    "E13" = 12.76635966044616
    "F13" = 36.96497799998596
    "I13" = 27.65037256314134
    "J13" = 10.08083020395006
    "K13" = 11.38984662116081
    "L13" = 11.17585043133108
    "M13" = 15.52867732109573
    "N13" = 21.15957168292993
    "O13" = 28.46807876366548
    "B14" = 14.21695584116649
    "C14" = 6.973969392220628
    "D14" = 8.02514593675939
    "E14" = 12.76936770067835
    "F14" = 36.97420461566063
    "I14" = 27.66357079098093
    "J14" = 10.08343185157588
    "K14" = 11.36170633769756
    "L14" = 11.17609455799619
    "M14" = 15.52145630180591
    "N14" = 21.16985330876427
    "O14" = 28.47906323873795
    "B15" = 14.19313182553993
    "C15" = 6.966492332626085
    "D15" = 8.023109889071518
    "E15" = 12.77123324146964
    "F15" = 36.97998385345214
    "I15" = 27.67173285727839
    "J15" = 10.085035023851
    "K15" = 11.34446277734852
    "L15" = 11.17626631341258
    "M15" = 15.51706629370595
    "N15" = 21.17618294082077
    "O15" = 28.4858807295874
    "B16" = 14.05640628794962
    "C16" = 6.923322764841343
    "D16" = 8.011717248784466
    "E16" = 12.78227728262143
    "F16" = 37.01505920267733
    "I16" = 27.71969852637564
    "J16" = 10.09436934809808
    "K16" = 11.24556727913621
    "L16" = 11.17759185006993
    "M16" = 15.49241987841719
    "N16" = 21.21294430666192
    "O16" = 28.52631576793691
    "B17" = 13.97239125124857
    "C17" = 6.896557636723802
    "D17" = 8.00497577827978
    "E17" = 12.7893682104558
    "F17" = 37.03832293529368
    "I17" = 27.75018787189503
    "J17" = 10.1002272514695
    "K17" = 11.18485441608092
    "L17" = 11.17871029005581
    "M17" = 15.47776031418879
    "N17" = 21.23593324688768
    "O17" = 28.55234072949432
    "B18" = 13.92403315777436
    "C18" = 6.881061387108028
    "D18" = 8.001190032024505
    "E18" = 12.79356294045101
    "F18" = 37.05234528727024
    "I18" = 27.7681156065274
    "J18" = 10.10364502597965
    "K18" = 11.14992886947708
    "L18" = 11.17946605948127
    "M18" = 15.46949931984689
    "N18" = 21.24931682941508
    "O18" = 28.56775762758414
    "B19" = 13.90765547971437
    "C19" = 6.87579735297486
    "D19" = 7.999924079150644
    "E19" = 12.7950031804445
    "F19" = 37.05720321614396
    "I19" = 27.77425280245012
    "J19" = 10.10481055902015
    "K19" = 11.13810384608543
    "L19" = 11.17974129408724
    "M19" = 15.46673178010934
    "N19" = 21.25387595693292
    "O19" = 28.57305446480785
    "B20" = 13.98133876667077
    "C20" = 6.899417396177525
    "D20" = 8.005683939302688
    "E20" = 12.78860134506545
    "F20" = 37.03578005978272
    "I20" = 27.74690175631224
    "J20" = 10.09959865417165
    "K20" = 11.19131818117187
    "L20" = 11.17857959615918
    "M20" = 15.47930321000913
    "N20" = 21.23346938478401
    "O20" = 28.54952395654761
    "B21" = 14.22837601519998
    "C21" = 6.977549088394463
    "D21" = 8.026127114535168
    "E21" = 12.76847936803339
    "F21" = 36.97146792902772
    "I21" = 27.65967794456631
    "J21" = 10.082665687036
    "K21" = 11.36997330270619
    "L21" = 11.17601822295707
    "M21" = 15.52357040712993
    "N21" = 21.16682671262603
    "O21" = 28.47581824615075
    "B22" = 14.38918264992777
    "C22" = 7.027665811749707
    "D22" = 8.040287901761753
    "E22" = 12.75635969552677
    "F22" = 36.9351503643577
    "I22" = 27.60615315172103
    "J22" = 10.0720284225214
    "K22" = 11.48645862972199
    "L22" = 11.17533955055053
    "M22" = 15.55398647143798
    "N22" = 21.12469785150494
    "O22" = 28.43163868389415
    "B23" = 14.30343225216302
    "C23" = 7.001006193634685
    "D23" = 8.032657678492495
    "E23" = 12.76273401738338
    "F23" = 36.95401092713512
    "I23" = 27.63440189589185
    "J23" = 10.07766654349532
    "K23" = 11.42432449167555
    "L23" = 11.17561100710494
    "M23" = 15.53761876566588
    "N23" = 21.14705269092995
    "O23" = 28.45485286245034
    "B24" = 13.97729376786151
    "C24" = 6.898124836808926
    "D24" = 8.005363499140259
    "E24" = 12.78894767712714
    "F24" = 37.03692767620929
    "I24" = 27.7483861662821
    "J24" = 10.09988268727963
    "K24" = 11.18839597551659
    "L24" = 11.1786383315596
    "M24" = 15.47860514627685
    "N24" = 21.23458277676426
    "O24" = 28.55079600294292
    "B25" = 13.6243779010065
    "C25" = 6.783289944601287
    "D25" = 7.979439743450039
    "E25" = 12.82166866513225
    "F25" = 37.15089665406124
    "I25" = 27.88628909951379
    "J25" = 10.12570061094156
    "K25" = 10.93385541514496
    "L25" = 11.18620704500258
    "M25" = 15.42149823581666
    "N25" = 21.33515633216837
    "O25" = 28.67141377487469
}

foreach ($key in $values.Keys) {
    $ws.Range($key).Value = $values[$key]
}

